$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new header cell A1 with value "Accion"
$ws.Range("A1").Value = "Accion"

# Columns B:C best-fit width
$ws.Columns("B:C").AutoFit() | Out-Null

# Move selection to A2 (matches the saved selection in the diff)
$ws.Range("A2").Select() | Out-Null
